$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = 44062
$ws.Range("B37").Value = 2
$ws.Range("C37").Value = "User update toimimaan ja tokenien muuttelua"

$ws.Range("H33").Select()
